{"js": "// Replace each two-digit multiplication equation in the document with its\n// updated version, per the authoritative diff. Every \"old\" value below is\n// unique within the document, so an exact, case-sensitive whole-match\n// search/replace is unambiguous and safe.\nconst replacements = [\n  [\"73\u00d715=\", \"72\u00d792=\"],\n  [\"23\u00d711=\", \"69\u00d733=\"],\n  [\"28\u00d739=\", \"34\u00d734=\"],\n  [\"60\u00d736=\", \"72\u00d770=\"],\n  [\"67\u00d739=\", \"99\u00d739=\"],\n  [\"36\u00d783=\", \"75\u00d734=\"],\n  [\"93\u00d791=\", \"32\u00d720=\"],\n  [\"21\u00d712=\", \"65\u00d770=\"],\n  [\"37\u00d758=\", \"86\u00d769=\"],\n  [\"45\u00d733=\", \"51\u00d729=\"],\n  [\"79\u00d797=\", \"65\u00d749=\"],\n  [\"24\u00d744=\", \"52\u00d769=\"],\n  [\"34\u00d776=\", \"87\u00d735=\"],\n  [\"63\u00d771=\", \"90\u00d749=\"],\n  [\"92\u00d738=\", \"63\u00d742=\"],\n  [\"16\u00d717=\", \"74\u00d746=\"],\n  [\"89\u00d786=\", \"36\u00d734=\"],\n  [\"25\u00d724=\", \"27\u00d735=\"],\n  [\"95\u00d723=\", \"79\u00d752=\"],\n  [\"14\u00d722=\", \"52\u00d770=\"],\n  [\"25\u00d761=\", \"56\u00d782=\"],\n  [\"51\u00d732=\", \"64\u00d715=\"],\n  [\"21\u00d725=\", \"47\u00d743=\"],\n  [\"26\u00d714=\", \"11\u00d719=\"],\n  [\"51\u00d780=\", \"92\u00d723=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the document with its\n# updated version, per the authoritative diff. Every \"old\" value below is\n# unique within the document, so an exact, case-sensitive whole-text\n# Find/Replace is unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"73\u00d715=\"; New = \"72\u00d792=\" },\n    @{ Old = \"23\u00d711=\"; New = \"69\u00d733=\" },\n    @{ Old = \"28\u00d739=\"; New = \"34\u00d734=\" },\n    @{ Old = \"60\u00d736=\"; New = \"72\u00d770=\" },\n    @{ Old = \"67\u00d739=\"; New = \"99\u00d739=\" },\n    @{ Old = \"36\u00d783=\"; New = \"75\u00d734=\" },\n    @{ Old = \"93\u00d791=\"; New = \"32\u00d720=\" },\n    @{ Old = \"21\u00d712=\"; New = \"65\u00d770=\" },\n    @{ Old = \"37\u00d758=\"; New = \"86\u00d769=\" },\n    @{ Old = \"45\u00d733=\"; New = \"51\u00d729=\" },\n    @{ Old = \"79\u00d797=\"; New = \"65\u00d749=\" },\n    @{ Old = \"24\u00d744=\"; New = \"52\u00d769=\" },\n    @{ Old = \"34\u00d776=\"; New = \"87\u00d735=\" },\n    @{ Old = \"63\u00d771=\"; New = \"90\u00d749=\" },\n    @{ Old = \"92\u00d738=\"; New = \"63\u00d742=\" },\n    @{ Old = \"16\u00d717=\"; New = \"74\u00d746=\" },\n    @{ Old = \"89\u00d786=\"; New = \"36\u00d734=\" },\n    @{ Old = \"25\u00d724=\"; New = \"27\u00d735=\" },\n    @{ Old = \"95\u00d723=\"; New = \"79\u00d752=\" },\n    @{ Old = \"14\u00d722=\"; New = \"52\u00d770=\" },\n    @{ Old = \"25\u00d761=\"; New = \"56\u00d782=\" },\n    @{ Old = \"51\u00d732=\"; New = \"64\u00d715=\" },\n    @{ Old = \"21\u00d725=\"; New = \"47\u00d743=\" },\n    @{ Old = \"26\u00d714=\"; New = \"11\u00d719=\" },\n    @{ Old = \"51\u00d780=\"; New = \"92\u00d723=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
